# Apply the latest cryptos list update (prices, 1h volume %, and rank reshuffles)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.216.91"
$ws.Range("E2").Value = "  -0.13%  "
# Row 3
$ws.Range("D3").Value = "1.855.02"
$ws.Range("E3").Value = "  -0.49%  "
# Row 4
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.10%  "
# Row 5
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").Value = "'241.33"
$ws.Range("E5").Value = "  -0.49%  "
# Row 6
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").Value = "'0.6997"
$ws.Range("E6").Value = "  -0.80%  "
# Row 7
$ws.Range("E7").Value = "  -0.04%  "
# Row 8
$ws.Range("D8").Value = "'0.3090"
$ws.Range("E8").Value = "  -0.61%  "
# Row 9
$ws.Range("D9").Value = "'0.07710"
$ws.Range("E9").Value = "  -1.49%  "
# Row 10
$ws.Range("D10").Value = "'23.78"
$ws.Range("E10").Value = "  -2.01%  "
# Row 11
$ws.Range("D11").Value = "'0.07797"
$ws.Range("E11").Value = "  -2.43%  "
# Row 12
$ws.Range("D12").Value = "1.860.76"
$ws.Range("E12").Value = "  -2.03%  "
# Row 13
$ws.Range("D13").Value = "'92.13"
$ws.Range("E13").Value = "  -1.54%  "
# Row 14
$ws.Range("D14").Value = "'5.093"
$ws.Range("E14").Value = "  -1.69%  "
# Row 15
$ws.Range("D15").Value = "'0.6870"
$ws.Range("E15").Value = "  -1.10%  "
# Row 16
$ws.Range("D16").Value = "'6.493"
$ws.Range("E16").Value = "  +2.10%  "
# Row 17
$ws.Range("D17").Value = "'0.000008380"
$ws.Range("E17").Value = "  +1.23%  "
# Row 18
$ws.Range("D18").Value = "29.231.72"
$ws.Range("E18").Value = "  -1.01%  "
# Row 19
$ws.Range("D19").Value = "'249.38"
$ws.Range("E19").Value = "  -0.82%  "
# Row 20
$ws.Range("D20").Value = "2.115.03"
$ws.Range("E20").Value = "  -2.80%  "
# Row 21
$ws.Range("D21").Value = "'12.84"
$ws.Range("E21").Value = "  -2.09%  "
# Row 22
$ws.Range("E22").Value = "  -0.10%  "
# Row 23
$ws.Range("D23").Value = "'7.535"
$ws.Range("E23").Value = "  +0.09%  "
# Row 24
$ws.Range("E24").Value = "  -0.06%  "
# Row 25
$ws.Range("D25").Value = "'0.1524"
$ws.Range("E25").Value = "  -1.96%  "
# Row 26
$ws.Range("D26").Value = "'160.13"
$ws.Range("E26").Value = "  +0.29%  "
# Row 27
$ws.Range("D27").Value = "'8.841"
$ws.Range("E27").Value = "  -1.68%  "
# Row 28
$ws.Range("D28").Value = "'18.50"
$ws.Range("E28").Value = "  -0.94%  "
# Row 29
$ws.Range("E29").Value = "  +4.07%  "
# Row 30
$ws.Range("D30").Value = "'4.233"
$ws.Range("E30").Value = "  -0.94%  "
# Row 31
$ws.Range("D31").Value = "'4.200"
$ws.Range("E31").Value = "  -1.70%  "
# Row 32
$ws.Range("E32").Value = "  -1.74%  "
# Row 33
$ws.Range("D33").Value = "'0.05187"
# Row 34
$ws.Range("D34").Value = "'0.7606"
$ws.Range("E34").Value = "  +1.99%  "
# Row 35
$ws.Range("D35").Value = "'1.841"
$ws.Range("E35").Value = "  -2.68%  "
# Row 36
$ws.Range("D36").Value = "'1.161"
$ws.Range("E36").Value = "  +0.52%  "
# Row 37
$ws.Range("D37").Value = "'2.712"
$ws.Range("E37").Value = "  +0.31%  "
# Row 38
$ws.Range("D38").Value = "'0.01860"
$ws.Range("E38").Value = "  +0.07%  "
# Row 39
$ws.Range("D39").Value = "1.222.66"
$ws.Range("E39").Value = "  -1.46%  "
# Row 40
$ws.Range("D40").Value = "'2.724"
$ws.Range("E40").Value = "  -0.77%  "
# Row 41
$ws.Range("D41").Value = "'0.8958"
$ws.Range("E41").Value = "  -0.70%  "
# Row 42
$ws.Range("D42").Value = "'109.72"
$ws.Range("E42").Value = "  -1.24%  "
# Row 43
$ws.Range("E43").Value = "  -0.01%  "
# Row 44
$ws.Range("D44").Value = "'5.552"
$ws.Range("E44").Value = "  -10.89%  "
# Row 45
$ws.Range("D45").Value = "2.015.20"
$ws.Range("E45").Value = "  -2.84%  "
# Row 46
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "'0.00000000124"
$ws.Range("E46").Value = "  +1.20%  "
# Row 47
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'65.09"
$ws.Range("E47").Value = "  -9.35%  "
# Row 48
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value = "'0.5183"
$ws.Range("E48").Value = "  -0.33%  "
# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.528"
$ws.Range("E49").Value = "  +1.53%  "
# Row 50
$ws.Range("D50").Value = "'1.747"
$ws.Range("E50").Value = "  -2.36%  "
# Row 51
$ws.Range("D51").Value = "'6.998"
$ws.Range("E51").Value = "  +0.38%  "
